$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.891649961471558
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.942132234573364
$ws.Range("D1").Value = 2.909365653991699
$ws.Range("E1").Value = 2.460869073867798
